# LOM3233.xlsx content refresh:
#  - "Objetivos" cell now holds the responsible-professor string instead of
#    the long objectives text
#  - "Docentes responsáveis" text row was removed, shifting everything below
#    it up by one row
#  - "Programa" full-syllabus text was removed (another row collapses)
#  - "Critério"/"Norma de recuperação"/"Bibliografia" labels shift up one row
#    each and the old bibliography paragraph is dropped entirely
#  - the final row (old row 25, "LOM3221 ..." requisite) is removed since
#    everything below "Requisitos:" shifts up by one row too

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 10: Objetivos -------------------------------------------------
$ws.Range("B10").Value = "519033 - Carlos Yujiro Shigue"
$ws.Range("C10").Value = "519033 - Carlos Yujiro Shigue"

# --- Row 13: now "Programa resumido:" / "Semestral" ---------------------
$ws.Range("A13").Value = "Programa resumido:"
$ws.Range("B13").Value = "Semestral"
$ws.Range("C13").Value = "Semestral"
$ws.Rows.Item(13).RowHeight = 60

# --- Row 14: now just "Short syllabus:" ----------------------------------
$ws.Range("A14").Value = "Short syllabus:"
$ws.Range("B14").ClearContents()
$ws.Range("C14").ClearContents()

# --- Row 15: now "Programa:" / "01/01/2012" ------------------------------
$ws.Range("A15").Value = "Programa:"
$ws.Range("B15").Value = "01/01/2012"
$ws.Range("C15").Value = "01/01/2012"
$ws.Rows.Item(15).RowHeight = 120

# --- Row 16: now just "Syllabus:" ----------------------------------------
$ws.Range("A16").Value = "Syllabus:"
$ws.Range("B16").ClearContents()
$ws.Range("C16").ClearContents()

# --- Row 17: now just "Avaliação:" ---------------------------------------
$ws.Range("A17").Value = "Avaliação:"
$ws.Rows.Item(17).RowHeight = 15

# --- Row 18: now "Método:" / docente string ------------------------------
$ws.Range("A18").Value = "Método:"
$ws.Range("B18").Value = "519033 - Carlos Yujiro Shigue"
$ws.Range("C18").Value = "519033 - Carlos Yujiro Shigue"
$ws.Rows.Item(18).RowHeight = 60

# --- Row 19: now "Critério:" ----------------------------------------------
$ws.Range("A19").Value = "Critério:"

# --- Row 20: now "Norma de recuperação:" ----------------------------------
$ws.Range("A20").Value = "Norma de recuperação:"

# --- Row 21: now "Bibliografia:" / the "Aplicação..." paragraph ----------
$ws.Range("A21").Value = "Bibliografia:"
$ws.Range("B21").Value = "Aplicação de uma prova escrita dentro do prazo regimental antes do início do próximo semestre letivo. A nota da segunda avaliação será a média aritmética entre a nota da prova de recuperação e a nota final da primeira avaliação"
$ws.Range("C21").Value = "Aplicação de uma prova escrita dentro do prazo regimental antes do início do próximo semestre letivo. A nota da segunda avaliação será a média aritmética entre a nota da prova de recuperação e a nota final da primeira avaliação"
$ws.Rows.Item(21).RowHeight = 120

# --- Row 22: now just "Requisitos:" (old Bibliografia body text dropped) -
$ws.Range("A22").Value = "Requisitos:"
$ws.Range("B22").ClearContents()
$ws.Range("C22").ClearContents()
$ws.Rows.Item(22).RowHeight = 15

# --- Row 23: now the first requisite (A is empty) -------------------------
$ws.Range("A23").ClearContents()
$ws.Range("B23").Value = "LOM3206 -  Eletrônica  (Requisito)`n"
$ws.Range("C23").Value = "LOM3206 -  Eletrônica  (Requisito)`n"
$ws.Rows.Item(23).RowHeight = 30

# --- Row 24: now the second requisite -------------------------------------
$ws.Range("B24").Value = "LOM3221 -  Laboratório de Eletrônica  (Requisito)`n"
$ws.Range("C24").Value = "LOM3221 -  Laboratório de Eletrônica  (Requisito)`n"
$ws.Rows.Item(24).RowHeight = 30

# --- Old row 25 is no longer needed; remove it so the sheet collapses to
#     A1:C24 exactly as in the edited workbook.
$ws.Rows.Item(25).Delete()
